# Update crypto price/volume data to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.872.03"
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").Value = "'1.714.09"
$ws.Range("E3").Value = "  +0.82%  "

$ws.Range("D4").Value = "'0.9986"
$ws.Range("E4").Value = "  -0.47%  "

$ws.Range("D5").Value = "'318.05"
$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("D6").Value = "'0.9979"
$ws.Range("E6").Value = "  -0.46%  "

$ws.Range("D7").Value = "'0.3934"
$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("E8").Value = "  +0.20%  "

$ws.Range("E9").Value = "  -0.66%  "

$ws.Range("D10").Value = "'53.67"

$ws.Range("D11").Value = "'0.9966"
$ws.Range("E11").Value = "  -0.76%  "

$ws.Range("D12").Value = "'0.08844"
$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("D13").Value = "'26.33"
$ws.Range("E13").Value = "  +10.85%  "

$ws.Range("D14").Value = "'7.536"
$ws.Range("E14").Value = "  -1.65%  "

$ws.Range("D15").Value = "'8.133"
$ws.Range("E15").Value = "  -0.51%  "

$ws.Range("E16").Value = "  +2.85%  "

$ws.Range("D17").Value = "'1.718.23"
$ws.Range("E17").Value = "  +1.57%  "

$ws.Range("D18").Value = "'97.43"
$ws.Range("E18").Value = "  -2.19%  "

$ws.Range("D19").Value = "'0.07215"
$ws.Range("E19").Value = "  +1.43%  "

$ws.Range("D20").Value = "'20.84"
$ws.Range("E20").Value = "  +4.81%  "

$ws.Range("E21").Value = "  +2.78%  "

$ws.Range("D22").Value = "'0.9981"
$ws.Range("E22").Value = "  -0.69%  "

$ws.Range("D23").Value = "'14.47"
$ws.Range("E23").Value = "  -1.79%  "

$ws.Range("D24").Value = "'24.873.47"
$ws.Range("E24").Value = "  +0.56%  "

$ws.Range("D25").Value = "'3.026"
$ws.Range("E25").Value = "  -4.23%  "

$ws.Range("D26").Value = "'2.335"
$ws.Range("E26").Value = "  -0.70%  "

$ws.Range("D27").Value = "'23.58"
$ws.Range("E27").Value = "  +2.18%  "

$ws.Range("D28").Value = "'167.10"
$ws.Range("E28").Value = "  +1.60%  "

$ws.Range("D29").Value = "'5.975"
$ws.Range("E29").Value = "  +16.16%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'146.47"
$ws.Range("E30").Value = "  +5.21%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'8.557"
$ws.Range("E31").Value = "  -6.84%  "

$ws.Range("B32").Value = "WEMIXTOKEN"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").Value = "'2.261"
$ws.Range("E32").Value = "  +15.41%  "

$ws.Range("B33").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C33").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D33").Value = "'1.905.66"
$ws.Range("E33").Value = "  +1.52%  "

$ws.Range("D34").Value = "'0.08875"
$ws.Range("E34").Value = "  -1.74%  "

$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "'0.03180"
$ws.Range("E35").Value = "  +5.05%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.061"
$ws.Range("E36").Value = "  -1.02%  "

$ws.Range("D37").Value = "'7.298"
$ws.Range("E37").Value = "  -8.24%  "

$ws.Range("D38").Value = "'0.2863"
$ws.Range("E38").Value = "  +2.73%  "

$ws.Range("D39").Value = "'0.8490"
$ws.Range("E39").Value = "  +9.27%  "

$ws.Range("D40").Value = "'10.94"
$ws.Range("E40").Value = "  -1.68%  "

$ws.Range("D41").Value = "'0.09271"
$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("E42").Value = "  -1.69%  "

$ws.Range("E43").Value = "  +1.39%  "

$ws.Range("D44").Value = "'17.50"
$ws.Range("E44").Value = "  +8.15%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'2.712"
$ws.Range("E45").Value = "  +3.05%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.7488"
$ws.Range("E46").Value = "  +3.51%  "

$ws.Range("D47").Value = "'4.277"
$ws.Range("E47").Value = "  +1.23%  "

$ws.Range("D48").Value = "'1.407"
$ws.Range("E48").Value = "  +3.32%  "

$ws.Range("D49").Value = "'0.9983"
$ws.Range("E49").Value = "  -0.35%  "

$ws.Range("D50").Value = "'140.83"
$ws.Range("E50").Value = "  +0.38%  "

$ws.Range("D51").Value = "'0.08304"
$ws.Range("E51").Value = "  +3.87%  "
